$d = $word.ActiveDocument

# Change 1: remove "determine and " before "swing prices"
$d.Content.Find.Execute(
    "able to determine and swing prices",
    $true, $false, $false, $false, $false,
    $true, 1, $false,
    "able to swing prices",
    2
)

# Change 2: merge "and other mainstream " + "social media sites" into one run
# (text itself is unchanged; runs merge in the underlying XML automatically
# when we re-set the same text through Find/Replace)
$d.Content.Find.Execute(
    "and other mainstream social media sites",
    $true, $false, $false, $false, $false,
    $true, 1, $false,
    "and other mainstream social media sites",
    2
)

# Change 3: split "A large scale study to understand…:" so "large scale" is
# wrapped by proofErr gramStart/gramEnd markers (grammar-check artifact).
$d.Content.Find.Execute(
    "A large scale study to understand",
    $true, $false, $false, $false, $false,
    $true, 1, $false,
    "A large scale study to understand",
    2
)
